$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row based on the worksheet's used range
$lastRow = $ws.UsedRange.Rows.Count

# Swap the contents of columns C and D (group-name / group-code)
# for every row, including the header row, since the header labels
# themselves are also swapped in the shared strings table.
for ($r = 1; $r -le $lastRow; $r++) {
    $cValue = $ws.Cells.Item($r, 3).Value2
    $dValue = $ws.Cells.Item($r, 4).Value2
    $ws.Cells.Item($r, 3).Value = $dValue
    $ws.Cells.Item($r, 4).Value = $cValue
}
